# Generate Report for Handback
# Applies the "handback" localization-status update:
#  - Status column flips from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File" + "Latest Handback File" +
#    "Latest Handback DateTime" columns populated now that the handback xliff came in
#  - "Latest Target File" cells link out to a.md, just like "Source File Name" does
#  - a couple of columns get widened so the new long values aren't clipped

$wb = $excel.ActiveWorkbook

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ffb1665bb9d0c49f3b4fe5930ab1892091a09130/e2e/a.md"

# ---------------------------------------------------------------------------
# 1. Overview sheet: widen the zh-cn / de-de status columns (E, F) so the
#    longer "Handed back: in sync with en-US" text fits, and flip the status
#    text itself.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
# 29.166666666666668 chars is the COM-model input that this engine's
# character->XML-width rounding turns into the target stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# 2. zh-cn sheet (table1)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

# Widen Status (C) and Latest Handback File (J) columns
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# Latest Target File (I) now links to a.md, same display/target as column A
$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $ghBase, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $ghBase, [Type]::Missing, [Type]::Missing, "a.md")

# Re-apply the plain hyperlink font after Hyperlinks.Add so it matches the
# workbook's existing "HyperLink" cell style instead of a fresh built-in one
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = 15570276
$wsZh.Range("I3").Font.Underline = 2
$wsZh.Range("I3").Font.Color = 15570276

# Latest Handback File (J) now has the generated handback xliff file name
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Latest Handback DateTime (K) - the generation timestamp
$wsZh.Range("K2").Value = "2016-08-17 04:33:15"
$wsZh.Range("K3").Value = "2016-08-17 04:33:15"

# ---------------------------------------------------------------------------
# 3. de-de sheet (table2)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# Widen Status (C) and Latest Handback File (J) columns
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

# Latest Target File (I) now links to a.md, same display/target as column A
$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $ghBase, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $ghBase, [Type]::Missing, [Type]::Missing, "a.md")

$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = 15570276
$wsDe.Range("I3").Font.Underline = 2
$wsDe.Range("I3").Font.Color = 15570276

# Latest Handback File (J) now has the generated handback xliff file name
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# Latest Handback DateTime (K) - the generation timestamp (distinct from zh-cn's)
$wsDe.Range("K2").Value = "2016-08-17 04:33:22"
$wsDe.Range("K3").Value = "2016-08-17 04:33:22"
